$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: set directly ---
$ws.Range("D2").Value = "60.179.91"
$ws.Range("E2").Value = "  -4.32%  "
$ws.Range("D3").Value = "3.300.41"
$ws.Range("E3").Value = "  -4.40%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -4.10%  "
$ws.Range("E6").Value = "  -4.91%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.302.59"
$ws.Range("E8").Value = "  -4.34%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("E11").Value = "  -4.46%  "
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").Value = "3.855.18"
$ws.Range("E13").Value = "  -4.66%  "
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("D16").Value = "3.297.93"
$ws.Range("E16").Value = "  -4.41%  "
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("D18").Value = "60.000.33"
$ws.Range("E18").Value = "  -4.61%  "
$ws.Range("E19").Value = "  -5.07%  "
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("E21").Value = "  -5.22%  "
$ws.Range("E22").Value = "  -3.52%  "
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "3.438.02"
$ws.Range("E26").Value = "  -4.24%  "
$ws.Range("E27").Value = "  -10.65%  "
$ws.Range("E28").Value = "  -7.09%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  -7.51%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -4.90%  "
$ws.Range("E33").Value = "  -5.46%  "
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("E35").Value = "  -8.32%  "
$ws.Range("E36").Value = "  -4.77%  "
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("E39").Value = "  -8.06%  "
$ws.Range("E40").Value = "  -16.43%  "
$ws.Range("D41").Value = "3.323.78"
$ws.Range("E41").Value = "  -4.69%  "
$ws.Range("E42").Value = "  -5.13%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  -4.80%  "
$ws.Range("E45").Value = "  -4.78%  "
$ws.Range("E46").Value = "  -6.53%  "
$ws.Range("E47").Value = "  -5.81%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.348.85"
$ws.Range("E48").Value = "  -8.30%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("E50").Value = "  -6.15%  "
$ws.Range("E51").Value = "  -7.68%  "

# --- Numeric-looking text values: force text format to preserve exact string, then restore default style ---
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.119"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.406"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000166"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.544"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000103"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0740"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.750"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.17"
$ws.Range("D51").Style = "Normal"
